# Update PLC data 2025-10-13 14:13:35
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7205
$ws.Range("C3").Value = 177214
$ws.Range("C4").Value = 167170
$ws.Range("C8").Value = 64.72
